$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, pushing existing rows 47:169 down to 48:170.
$ws.Rows("47:47").Insert()

# Populate the newly inserted row 47 with the new weekly data point.
$ws.Cells.Item(47, 1).Value = 8
$ws.Cells.Item(47, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(47, 3).Value = "Coquimbo"
$ws.Cells.Item(47, 4).Value = 45148
$ws.Cells.Item(47, 5).Value = 4
$ws.Cells.Item(47, 6).Value = 100114007
$ws.Cells.Item(47, 7).Value = "Jengibre"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 360
$ws.Cells.Item(47, 11).Value = 18000
$ws.Cells.Item(47, 12).Value = 19000
$ws.Cells.Item(47, 13).Value = 18500
$ws.Cells.Item(47, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(47, 15).Value = "Perú"
$ws.Cells.Item(47, 16).Value = 1423
$ws.Cells.Item(47, 17).Value = 13
$ws.Cells.Item(47, 18).Value = "Hortaliza"
